$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new label cell next to the "Member list" title, referencing a new shared string.
$ws.Range("C1").Value = "(in sample.xlsx)"

# Update column widths: B:F all get the same width, replacing the previous
# per-column custom widths (B=5.71, D=14.43 bestFit, E=6).
$ws.Range("B1:F1").ColumnWidth = 13.5708705357

# Move/normalize the saved selection to C1 (matches the new content focus).
[void]$ws.Range("C1").Select()
